$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: fill in previously-empty Q:W cells ---
$ws.Range("Q16").Value = "3h 19m"
$ws.Range("R16").Value = "3h 12m"
$ws.Range("S16").Value = 1.58233709667704
$ws.Range("T16").Value = 1.33264932588068
$ws.Range("U16").Value = 1.95503590731588
$ws.Range("V16").Value = 1.66707482368021
$ws.Range("W16").Value = 0.0249292056730672

# --- Row 17: fill in previously-empty Q:W cells ---
$ws.Range("Q17").Value = "2h 36m"
$ws.Range("R17").Value = "1h 57m"
$ws.Range("S17").Value = 0.932814737119131
$ws.Range("T17").Value = 0.789893631753853
$ws.Range("U17").Value = 1.12210443679312
$ws.Range("V17").Value = 0.974634833775511
$ws.Range("W17").Value = 0.0151928268355334

# --- Update the selected cell shown when the workbook is reopened ---
$ws.Range("L20").Select()
